$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Copy() | Out-Null
$ws.Range("A68").PasteSpecial(-4122) | Out-Null

$ws.Range("C67").Copy() | Out-Null
$ws.Range("C68").PasteSpecial(-4122) | Out-Null

$ws.Range("B68").Value = "champignon"
$ws.Range("C68").Value = "WNP"
$ws.Range("A68").Value = "4de94b55e-538e-4225-93f3-303390e81ed8"

$ws.Range("A69").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1
